{"js": "// Office.js (Word JavaScript API) script.\n// Applies two textual edits:\n//  1. \"August 21\" + \" 2023\" (previously split across two runs) -> a single\n//     run reading \"August 21 2023\". Office.js `body.search()` can locate\n//     text that spans a run boundary, and rewriting that found range with\n//     `insertText(..., Replace)` collapses it into one run.\n//  2. \"control over our lives\" -> \"control our lives\" (drop stray \"over\").\n\nconst body = context.document.body;\n\n// --- Edit 1: normalize the date line -------------------------------------\nconst dateResults = body.search(\"August 21 2023\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"August 21 2023\", Word.InsertLocation.replace);\n} else {\n  // Defensive fallback in case the two runs aren't reported as one\n  // contiguous match: rewrite \"August 21\" and drop the old \" 2023\" tail.\n  const partial = body.search(\"August 21\", { matchCase: true });\n  partial.load(\"items\");\n  await context.sync();\n  if (partial.items.length > 0) {\n    partial.items[0].insertText(\"August 21 2023\", Word.InsertLocation.replace);\n    await context.sync();\n    const stray = body.search(\"2023 2023\", { matchCase: true });\n    stray.load(\"items\");\n    await context.sync();\n    if (stray.items.length > 0) {\n      stray.items[0].insertText(\"2023\", Word.InsertLocation.replace);\n    }\n  }\n}\n\nawait context.sync();\n\n// --- Edit 2: fix wording in the Sam Prentice paragraph --------------------\nconst wordingResults = body.search(\"control over our lives\", { matchCase: true });\nwordingResults.load(\"items\");\nawait context.sync();\n\nif (wordingResults.items.length > 0) {\n  wordingResults.items[0].insertText(\"control our lives\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies two textual edits:\n#  1. Merge the split \"August 21\" / \" 2023\" runs into a single run reading\n#     \"August 21 2023\" (Find.Execute can match text that spans a run\n#     boundary; replacing it with identical text collapses it to one run).\n#  2. Remove the stray \"over\" in \"control over our lives\" -> \"control our lives\".\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: normalize the date line -------------------------------------\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"August 21 2023\"\n$find1.Replacement.Text = \"August 21 2023\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2) | Out-Null\n\n# --- Edit 2: fix wording in the Sam Prentice paragraph --------------------\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"control over our lives\"\n$find2.Replacement.Text = \"control our lives\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n\n$d.Save()\n"}
